$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need the Absent (column H) value computed/filled based on
# the Real (column E) value: Absent = 1 when not Real, 0 when Real.
$rows = 4..15

foreach ($r in $rows) {
    $real = $ws.Cells.Item($r, 5).Value2
    if ($real -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
